$wb = $excel.ActiveWorkbook

# New values for column F (想去人数 / "want to go" count), rows 2-8
$newValues = @{
    2 = 2174
    3 = 1651
    4 = 323
    5 = 1064
    6 = 633
    7 = 33
    8 = 5748
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
